$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new Start/End numeric values for rows 3 and 4 (columns D/E)
$ws.Range("D3").Value = 70522
$ws.Range("E3").Value = 286126
$ws.Range("D4").Value = 107666
$ws.Range("E4").Value = 289620

# Update the selection to cover the header + data rows (A1:E4)
$ws.Range("A1:E4").Select()

$wb.Save()
